$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '68.377.90'
$ws.Cells.Item(2, 5).Value = '  +1.62%  '

$ws.Cells.Item(3, 4).Value = '3.909.27'
$ws.Cells.Item(3, 5).Value = '  +1.00%  '

$ws.Cells.Item(4, 4).Value = '0.999'
$ws.Cells.Item(4, 5).Value = '  -0.12%  '

$ws.Cells.Item(5, 4).Value = '480.46'
$ws.Cells.Item(5, 5).Value = '  +2.59%  '

$ws.Cells.Item(6, 4).Value = '144.84'
$ws.Cells.Item(6, 5).Value = '  -0.09%  '

$ws.Cells.Item(7, 4).Value = '0.621'
$ws.Cells.Item(7, 5).Value = '  -1.86%  '

$ws.Cells.Item(8, 4).Value = '0.997'
$ws.Cells.Item(8, 5).Value = '  -0.13%  '

$ws.Cells.Item(9, 4).Value = '0.725'
$ws.Cells.Item(9, 5).Value = '  -2.78%  '

$ws.Cells.Item(10, 4).Value = '0.167'
$ws.Cells.Item(10, 5).Value = '  +7.58%  '

$ws.Cells.Item(11, 4).Value = '0.0000352'
$ws.Cells.Item(11, 5).Value = '  +13.02%  '

$ws.Cells.Item(12, 4).Value = '42.65'
$ws.Cells.Item(12, 5).Value = '  -1.75%  '

$ws.Cells.Item(13, 4).Value = '10.64'
$ws.Cells.Item(13, 5).Value = '  +1.80%  '

$ws.Cells.Item(14, 4).Value = '4.530.36'
$ws.Cells.Item(14, 5).Value = '  +0.88%  '

$ws.Cells.Item(15, 4).Value = '14.63'
$ws.Cells.Item(15, 5).Value = '  -1.36%  '

$ws.Cells.Item(16, 4).Value = '3.937.88'
$ws.Cells.Item(16, 5).Value = '  +2.15%  '

$ws.Cells.Item(17, 5).Value = '  -0.40%  '

$ws.Cells.Item(18, 4).Value = '19.72'
$ws.Cells.Item(18, 5).Value = '  -1.71%  '

$ws.Cells.Item(19, 5).Value = '  -3.18%  '

$ws.Cells.Item(20, 4).Value = '68.371.05'
$ws.Cells.Item(20, 5).Value = '  +1.25%  '

$ws.Cells.Item(21, 4).Value = '435.66'
$ws.Cells.Item(21, 5).Value = '  -0.09%  '

$ws.Cells.Item(22, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(22, 4).Value = '14.66'
$ws.Cells.Item(22, 5).Value = '  -1.66%  '

$ws.Cells.Item(23, 2).Value = 'ImmutableX'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(23, 4).Value = '3.37'
$ws.Cells.Item(23, 5).Value = '  +1.60%  '

$ws.Cells.Item(24, 4).Value = '87.84'
$ws.Cells.Item(24, 5).Value = '  -1.40%  '

$ws.Cells.Item(25, 4).Value = '11.72'
$ws.Cells.Item(25, 5).Value = '  +18.08%  '

$ws.Cells.Item(26, 4).Value = '3.59'
$ws.Cells.Item(26, 5).Value = '  -0.14%  '

$ws.Cells.Item(27, 4).Value = '38.16'
$ws.Cells.Item(27, 5).Value = '  +0.29%  '

$ws.Cells.Item(28, 4).Value = '10.42'
$ws.Cells.Item(28, 5).Value = '  +2.91%  '

$ws.Cells.Item(29, 4).Value = '5.81'
$ws.Cells.Item(29, 5).Value = '  +4.72%  '

$ws.Cells.Item(30, 4).Value = '705.25'
$ws.Cells.Item(30, 5).Value = '  -3.29%  '

$ws.Cells.Item(31, 5).Value = '  -2.13%  '

$ws.Cells.Item(32, 4).Value = '13.36'
$ws.Cells.Item(32, 5).Value = '  -3.60%  '

$ws.Cells.Item(33, 5).Value = '  +2.64%  '

$ws.Cells.Item(34, 4).Value = '0.0₃0935'
$ws.Cells.Item(34, 5).Value = '  +37.81%  '

$ws.Cells.Item(35, 4).Value = '41.61'
$ws.Cells.Item(35, 5).Value = '  -6.03%  '

$ws.Cells.Item(36, 4).Value = '59.37'
$ws.Cells.Item(36, 5).Value = '  +1.77%  '

$ws.Cells.Item(37, 4).Value = '5.74'
$ws.Cells.Item(37, 5).Value = '  +4.33%  '

$ws.Cells.Item(38, 5).Value = '  -6.85%  '

$ws.Cells.Item(39, 4).Value = '0.999'
$ws.Cells.Item(39, 5).Value = '  -0.11%  '

$ws.Cells.Item(40, 5).Value = '  -2.18%  '

$ws.Cells.Item(41, 5).Value = '  +10.73%  '

$ws.Cells.Item(42, 4).Value = '2.75'
$ws.Cells.Item(42, 5).Value = '  +7.67%  '

$ws.Cells.Item(43, 5).Value = '  +2.54%  '

$ws.Cells.Item(44, 4).Value = '0.341'
$ws.Cells.Item(44, 5).Value = '  -1.73%  '

$ws.Cells.Item(45, 4).Value = '0.142'
$ws.Cells.Item(45, 5).Value = '  -0.27%  '

$ws.Cells.Item(46, 5).Value = '  -0.21%  '

$ws.Cells.Item(47, 5).Value = '  -0.86%  '

$ws.Cells.Item(48, 5).Value = '  -0.56%  '

$ws.Cells.Item(49, 4).Value = '146.02'
$ws.Cells.Item(49, 5).Value = '  +1.10%  '

$ws.Cells.Item(50, 4).Value = '3.14'
$ws.Cells.Item(50, 5).Value = '  -4.37%  '

$ws.Cells.Item(51, 4).Value = '2.85'
$ws.Cells.Item(51, 5).Value = '  -1.76%  '
